# Apply the strategy data update for 2024-02-22.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G values (rows 3-26)
$ws.Range("G3").Value = 53
$ws.Range("G4").Value = 51
$ws.Range("G5").Value = 51
$ws.Range("G6").Value = 51
$ws.Range("G7").Value = 51
$ws.Range("G8").Value = 51
$ws.Range("G9").Value = 51
$ws.Range("G10").Value = 51
$ws.Range("G11").Value = 51
$ws.Range("G12").Value = 51
$ws.Range("G13").Value = 51
$ws.Range("G14").Value = 51
$ws.Range("G15").Value = 51
$ws.Range("G16").Value = 51
$ws.Range("G17").Value = 51
$ws.Range("G18").Value = 51
$ws.Range("G19").Value = 55
$ws.Range("G20").Value = 51
$ws.Range("G21").Value = 51
$ws.Range("G22").Value = 51
$ws.Range("G23").Value = 51
$ws.Range("G24").Value = 51
$ws.Range("G25").Value = 51
$ws.Range("G26").Value = 51

# Update the view: scroll position (topLeftCell = A10) and active cell/selection (G20)
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G20").Select()
